$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44295
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 29000
$ws.Range("L2").Value = 30000
$ws.Range("M2").Value = 29500
$ws.Range("O2").Value = "Provincia del Elquí"
$ws.Range("P2").Value = 1180

$ws.Range("D3").Value = 44335
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 31000
$ws.Range("M3").Value = 30500
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 1220

$ws.Range("D4").Value = 44278
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 29000
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 29500
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 1180

$ws.Range("D5").Value = 44326
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 29000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = 29500
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 1180

$ws.Range("D6").Value = 44202
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 29000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 29500
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 1180

$ws.Range("D7").Value = 44330
$ws.Range("J7").Value = 520
$ws.Range("K7").Value = 32000
$ws.Range("L7").Value = 33000
$ws.Range("M7").Value = 32500
$ws.Range("O7").Value = "Provincia del Elquí"
$ws.Range("P7").Value = 1300

$ws.Range("D8").Value = 44302
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 28000
$ws.Range("L8").Value = 29000
$ws.Range("M8").Value = 28500
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 1140

$ws.Range("D9").Value = 44231
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 25000
$ws.Range("L9").Value = 26000
$ws.Range("M9").Value = 25500
$ws.Range("O9").Value = "Provincia del Elquí"
$ws.Range("P9").Value = 1020

$ws.Range("D10").Value = 44230
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 28000
$ws.Range("L10").Value = 29000
$ws.Range("M10").Value = 28500
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 1140

$ws.Range("D11").Value = 44272
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 28000
$ws.Range("L11").Value = 29000
$ws.Range("M11").Value = 28500
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 1140

$ws.Range("D12").Value = 44298
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 29000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 29500
$ws.Range("O12").Value = "Provincia del Elquí"
$ws.Range("P12").Value = 1180

$ws.Range("D13").Value = 44333
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 31500
$ws.Range("L13").Value = 32000
$ws.Range("M13").Value = 31750
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 1270

$ws.Range("D14").Value = 44293
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 29000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 29500
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 1180

$ws.Range("D15").Value = 44238
$ws.Range("J15").Value = 520
$ws.Range("K15").Value = 28000
$ws.Range("L15").Value = 29000
$ws.Range("M15").Value = 28500
$ws.Range("O15").Value = "Provincia del Elquí"
$ws.Range("P15").Value = 1140

$ws.Range("D16").Value = 44216
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 36000
$ws.Range("L16").Value = 37000
$ws.Range("M16").Value = 36500
$ws.Range("O16").Value = "Provincia del Elquí"
$ws.Range("P16").Value = 1460

$ws.Range("D17").Value = 44253
$ws.Range("J17").Value = 660
$ws.Range("K17").Value = 28000
$ws.Range("L17").Value = 29000
$ws.Range("M17").Value = 28500
$ws.Range("O17").Value = "Provincia del Elquí"
$ws.Range("P17").Value = 1140

$ws.Range("D18").Value = 44245
$ws.Range("J18").Value = 540
$ws.Range("K18").Value = 28000
$ws.Range("L18").Value = 29000
$ws.Range("M18").Value = 28500
$ws.Range("O18").Value = "Provincia del Elquí"
$ws.Range("P18").Value = 1140

$ws.Range("D19").Value = 44281
$ws.Range("J19").Value = 640
$ws.Range("K19").Value = 29000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 29500
$ws.Range("O19").Value = "Provincia del Elquí"
$ws.Range("P19").Value = 1180

$ws.Range("D20").Value = 44225
$ws.Range("J20").Value = 600
$ws.Range("K20").Value = 31000
$ws.Range("L20").Value = 32000
$ws.Range("M20").Value = 31500
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 1260

$ws.Range("D21").Value = 44321
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 29000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 29500
$ws.Range("O21").Value = "Provincia del Elquí"
$ws.Range("P21").Value = 1180

$ws.Range("D22").Value = 44300
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 29000
$ws.Range("L22").Value = 30000
$ws.Range("M22").Value = 29500
$ws.Range("O22").Value = "Provincia del Elquí"
$ws.Range("P22").Value = 1180

$ws.Range("D23").Value = 44218
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 34000
$ws.Range("L23").Value = 35000
$ws.Range("M23").Value = 34500
$ws.Range("O23").Value = "Provincia del Elquí"
$ws.Range("P23").Value = 1380

$ws.Range("D24").Value = 44223
$ws.Range("J24").Value = 660
$ws.Range("K24").Value = 32500
$ws.Range("L24").Value = 33000
$ws.Range("M24").Value = 32750
$ws.Range("O24").Value = "Provincia del Elquí"
$ws.Range("P24").Value = 1310

$ws.Range("D25").Value = 44246
$ws.Range("J25").Value = 600
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 29000
$ws.Range("M25").Value = 28500
$ws.Range("O25").Value = "Provincia del Elquí"
$ws.Range("P25").Value = 1140

$ws.Range("D26").Value = 44307
$ws.Range("J26").Value = 560
$ws.Range("K26").Value = 28000
$ws.Range("L26").Value = 29000
$ws.Range("M26").Value = 28500
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 1140

$ws.Range("D27").Value = 44274
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 28000
$ws.Range("L27").Value = 29000
$ws.Range("M27").Value = 28500
$ws.Range("O27").Value = "Provincia del Elquí"
$ws.Range("P27").Value = 1140

$ws.Range("D28").Value = 44309
$ws.Range("J28").Value = 600
$ws.Range("K28").Value = 28000
$ws.Range("L28").Value = 29000
$ws.Range("M28").Value = 28500
$ws.Range("O28").Value = "Provincia del Elquí"
$ws.Range("P28").Value = 1140

$ws.Range("D29").Value = 44200
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 31000
$ws.Range("L29").Value = 32000
$ws.Range("M29").Value = 31500
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 1260

$ws.Range("D30").Value = 44237
$ws.Range("J30").Value = 600
$ws.Range("K30").Value = 29000
$ws.Range("L30").Value = 30000
$ws.Range("M30").Value = 29500
$ws.Range("O30").Value = "Provincia del Elquí"
$ws.Range("P30").Value = 1180

$ws.Range("D31").Value = 44265
$ws.Range("J31").Value = 760
$ws.Range("K31").Value = 28000
$ws.Range("L31").Value = 29000
$ws.Range("M31").Value = 28500
$ws.Range("O31").Value = "Provincia del Elquí"
$ws.Range("P31").Value = 1140

$ws.Range("D32").Value = 44260
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 27000
$ws.Range("L32").Value = 28000
$ws.Range("M32").Value = 27500
$ws.Range("O32").Value = "Provincia del Elquí"
$ws.Range("P32").Value = 1100

$ws.Range("D33").Value = 44187
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 37000
$ws.Range("L33").Value = 38000
$ws.Range("M33").Value = 37500
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 1500

$ws.Range("D34").Value = 44264
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 28000
$ws.Range("L34").Value = 29000
$ws.Range("M34").Value = 28500
$ws.Range("O34").Value = "Provincia del Elquí"
$ws.Range("P34").Value = 1140

$ws.Range("D35").Value = 44251
$ws.Range("J35").Value = 700
$ws.Range("K35").Value = 29000
$ws.Range("L35").Value = 30000
$ws.Range("M35").Value = 29500
$ws.Range("O35").Value = "Provincia del Elquí"
$ws.Range("P35").Value = 1180

$ws.Range("D36").Value = 44252
$ws.Range("J36").Value = 520
$ws.Range("K36").Value = 28000
$ws.Range("L36").Value = 29000
$ws.Range("M36").Value = 28500
$ws.Range("O36").Value = "Provincia del Elquí"
$ws.Range("P36").Value = 1140

$ws.Range("D37").Value = 44221
$ws.Range("J37").Value = 460
$ws.Range("K37").Value = 35000
$ws.Range("L37").Value = 36000
$ws.Range("M37").Value = 35500
$ws.Range("O37").Value = "Provincia del Elquí"
$ws.Range("P37").Value = 1420

$ws.Range("D38").Value = 44316
$ws.Range("J38").Value = 600
$ws.Range("K38").Value = 28000
$ws.Range("L38").Value = 29000
$ws.Range("M38").Value = 28500
$ws.Range("O38").Value = "Provincia del Elquí"
$ws.Range("P38").Value = 1140

$ws.Range("D39").Value = 44279
$ws.Range("J39").Value = 560
$ws.Range("K39").Value = 29000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = 29500
$ws.Range("O39").Value = "Provincia del Elquí"
$ws.Range("P39").Value = 1180

$ws.Range("D40").Value = 44277
$ws.Range("J40").Value = 560
$ws.Range("K40").Value = 29000
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = 29500
$ws.Range("O40").Value = "Provincia del Elquí"
$ws.Range("P40").Value = 1180

$ws.Range("D41").Value = 44291
$ws.Range("J41").Value = 500
$ws.Range("K41").Value = 29000
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = 29500
$ws.Range("O41").Value = "Provincia del Elquí"
$ws.Range("P41").Value = 1180

$ws.Range("D42").Value = 44209
$ws.Range("J42").Value = 600
$ws.Range("K42").Value = 36000
$ws.Range("L42").Value = 37000
$ws.Range("M42").Value = 36500
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 1460

$ws.Range("D43").Value = 44244
$ws.Range("J43").Value = 640
$ws.Range("K43").Value = 29000
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = 29500
$ws.Range("O43").Value = "Provincia del Elquí"
$ws.Range("P43").Value = 1180

$ws.Range("D44").Value = 44239
$ws.Range("J44").Value = 600
$ws.Range("K44").Value = 28000
$ws.Range("L44").Value = 29000
$ws.Range("M44").Value = 28500
$ws.Range("O44").Value = "Provincia del Elquí"
$ws.Range("P44").Value = 1140

$ws.Range("D45").Value = 44323
$ws.Range("J45").Value = 600
$ws.Range("K45").Value = 29000
$ws.Range("L45").Value = 30000
$ws.Range("M45").Value = 29500
$ws.Range("O45").Value = "Provincia del Elquí"
$ws.Range("P45").Value = 1180
